$d = $word.ActiveDocument

# 1) "Đơn vị tổ chức" row: @QsUniversity -> @Unit
$d.Content.Find.Execute("@QsUniversity", $true, $false, $false, $false, $false,
                         $true, 1, $false, "@Unit", 2) | Out-Null

# 2) After @KeynoteSpeaker, append " thuộc trường đại học @QsUniversity"
#    as two additional runs in the same paragraph.
$found = $d.Content.Find.Execute("@KeynoteSpeaker", $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

$rng = $d.Content
$rng.Find.Execute("@KeynoteSpeaker") | Out-Null
$rng.Collapse(0)  # wdCollapseEnd

$rng.InsertAfter(" thuộc trường đại học ")
$rng.Collapse(0)
$rng.InsertAfter("@QsUniversity")
